# Update "Pre Experimental Phase" outcome measures (column B) on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value  = "A little worse "
$ws.Range("B8").Value  = "A little worse "
$ws.Range("B9").Value  = "A lot worse "
$ws.Range("B10").Value = "Somewhat worse "
$ws.Range("B12").Value = "A lot worse "
$ws.Range("B13").Value = "Somewhat worse "
$ws.Range("B15").Value = "A lot worse "
